# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values scraped from coinranking.com. Row 34/35 also swap which coin
# (Kaspa / ImmutableX) occupies that ranking slot, per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price values are forced to text (NumberFormat "@") so
# Excel doesn't coerce them into floating point numbers / lose trailing
# zeros, then the style is reset to "Normal" so no stray cell formatting
# is introduced.
$ws.Range("D2").Value = "67.124.19"
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("D3").Value = "3.537.61"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("D7").Value = "3.534.21"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "4.135.69"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "3.522.05"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D17").Value = "67.084.52"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "449.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.636"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "3.676.19"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.96%  "
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.159"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").Value = "3.533.40"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0870"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.890"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("E51").Value = "  -4.20%  "

Write-Output "Applied 90 cell changes"